# This workbook tracks weekly price observations for "Pepino dulce" at the
# Mapocho Venta Directa de Santiago market. Two new weekly rows need to be
# inserted right after the current row 124 (pushing the existing rows
# 125-132 down to 127-134), and populated with this week's data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows before the current row 125; everything currently at
# row 125 and below shifts down by two rows.
$ws.Rows("125:126").Insert()

# --- New row 125: "Primera" quality entry for the new date ---
$ws.Range("A125").Value2 = 12
$ws.Range("B125").Value2 = "Mapocho Venta Directa de Santiago"
$ws.Range("C125").Value2 = "Metropolitana"
$ws.Range("D125").Value2 = 44461
$ws.Range("E125").Value2 = 13
$ws.Range("F125").Value2 = 100112043
$ws.Range("G125").Value2 = "Pepino dulce"
$ws.Range("H125").Value2 = "Cultivar IV Región"
$ws.Range("I125").Value2 = "Primera"
$ws.Range("J125").Value2 = 510
$ws.Range("K125").Value2 = 20000
$ws.Range("L125").Value2 = 20000
$ws.Range("M125").Value2 = 20000
$ws.Range("N125").Value2 = "$/bandeja 18 kilos"
$ws.Range("O125").Value2 = "Provincia de Limarí"
$ws.Range("P125").Value2 = 1111
$ws.Range("Q125").Value2 = 18
$ws.Range("R125").Value2 = "Hortaliza"

# --- New row 126: "Segunda" quality entry for the new date ---
$ws.Range("A126").Value2 = 12
$ws.Range("B126").Value2 = "Mapocho Venta Directa de Santiago"
$ws.Range("C126").Value2 = "Metropolitana"
$ws.Range("D126").Value2 = 44461
$ws.Range("E126").Value2 = 13
$ws.Range("F126").Value2 = 100112043
$ws.Range("G126").Value2 = "Pepino dulce"
$ws.Range("H126").Value2 = "Cultivar IV Región"
$ws.Range("I126").Value2 = "Segunda"
$ws.Range("J126").Value2 = 450
$ws.Range("K126").Value2 = 17000
$ws.Range("L126").Value2 = 17000
$ws.Range("M126").Value2 = 17000
$ws.Range("N126").Value2 = "$/bandeja 18 kilos"
$ws.Range("O126").Value2 = "Provincia de Limarí"
$ws.Range("P126").Value2 = 944
$ws.Range("Q126").Value2 = 18
$ws.Range("R126").Value2 = "Hortaliza"
